$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Control 26
$ws.Range("D2").Value2 = 0.9999852330479696
$ws.Range("E2").Value2 = 0.9999852330479696

# Row 3 - Control 33
$ws.Range("C3").Value2 = $true
$ws.Range("D3").Value2 = 0.04570114962313795
$ws.Range("E3").Value2 = 0.04570114962313795

# Row 4 - Control 36
$ws.Range("D4").Value2 = 0.972847354826792
$ws.Range("E4").Value2 = 0.972847354826792

# Row 5 - Control 49
$ws.Range("D5").Value2 = 0.003132376245777385
$ws.Range("E5").Value2 = 0.003132376245777385

# Row 6 - Control 2
$ws.Range("D6").Value2 = 3.096311348770248 * [Math]::Pow(10, -10)
$ws.Range("E6").Value2 = 3.096311348770248 * [Math]::Pow(10, -10)

# Row 7 - MDD 36
$ws.Range("D7").Value2 = 0.9999999995284132
$ws.Range("E7").Value2 = 4.715867696347686 * [Math]::Pow(10, -10)

# Row 8 - MDD 10
$ws.Range("D8").Value2 = 0.9999904521862558
$ws.Range("E8").Value2 = 9.547813744181788 * [Math]::Pow(10, -6)

# Row 9 - MDD 39
$ws.Range("D9").Value2 = 0.9999999999965052
$ws.Range("E9").Value2 = 3.494760036915068 * [Math]::Pow(10, -12)

# Row 10 - MDD 14
$ws.Range("D10").Value2 = 0.9999278824782953
$ws.Range("E10").Value2 = 7.211752170466212 * [Math]::Pow(10, -5)

# Row 11 - MDD 18
$ws.Range("D11").Value2 = 0.999999999995999
$ws.Range("E11").Value2 = 4.001021736144139 * [Math]::Pow(10, -12)
$ws.Range("F11").Value2 = 1.47793984413147
$ws.Range("G11").Value2 = 0.8
